$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I20").Value = 0.4238415822648009
$ws.Range("J20").Value = 0.2208003802520098
$ws.Range("K20").Value = 0.2516512283222432
$ws.Range("L20").Value = 2.565679521489
